$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2067.4167
$ws.Range("I113").Value = 1843.3334
$ws.Range("J113").Value = 2291.5
$ws.Range("K113").Value = 1843.3334
$ws.Range("L113").Value = 2291.5
$ws.Range("M113").Value = 1410.6666
$ws.Range("N113").Value = -8799.5

$ws.Range("H137").Value = 2260.4673
$ws.Range("I137").Value = 921.875
$ws.Range("J137").Value = 3290.1538
$ws.Range("K137").Value = 2765.625
$ws.Range("L137").Value = 9870.4614
$ws.Range("M137").Value = -215.625
$ws.Range("N137").Value = -14970.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 51216.55
$ws.Range("I2").Value = 111972.336
$ws.Range("J2").Value = 1507.2727
$ws.Range("K2").Value = 111972.336
$ws.Range("L2").Value = 1507.2727
$ws.Range("M2").Value = -111859.336
$ws.Range("N2").Value = -1733.2727

$ws.Range("H32").Value = 3352.0876
$ws.Range("I32").Value = 3435.2642
$ws.Range("J32").Value = 2250
$ws.Range("K32").Value = 3435.2642
$ws.Range("L32").Value = 2250
$ws.Range("M32").Value = -3148.2642
$ws.Range("N32").Value = -2824

$ws.Range("H110").Value = 2651.9443
$ws.Range("I110").Value = 2785.6667
$ws.Range("J110").Value = 1983.3334
$ws.Range("K110").Value = 2785.6667
$ws.Range("L110").Value = 1983.3334
$ws.Range("M110").Value = -740.6667000000002
$ws.Range("N110").Value = -6073.3334

$ws.Range("H116").Value = 51216.55
$ws.Range("I116").Value = 111972.336
$ws.Range("J116").Value = 1507.2727
$ws.Range("K116").Value = 111972.336
$ws.Range("L116").Value = 1507.2727
$ws.Range("M116").Value = -109678.336
$ws.Range("N116").Value = -6095.2727

$ws.Range("H132").Value = 4608.877
$ws.Range("I132").Value = 2949.5
$ws.Range("K132").Value = 8848.5
$ws.Range("M132").Value = -6318.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 51216.55
$ws.Range("I3").Value = 111972.336
$ws.Range("J3").Value = 1507.2727
$ws.Range("K3").Value = 111972.336
$ws.Range("L3").Value = 1507.2727
$ws.Range("M3").Value = -111858.336
$ws.Range("N3").Value = -1735.2727

$ws.Range("H86").Value = 2527.0588
$ws.Range("I86").Value = 1512.5
$ws.Range("J86").Value = 3428.889
$ws.Range("K86").Value = 1512.5
$ws.Range("L86").Value = 3428.889
$ws.Range("M86").Value = -389.5
$ws.Range("N86").Value = -5674.889

$ws.Range("H89").Value = 2527.0588
$ws.Range("I89").Value = 1512.5
$ws.Range("J89").Value = 3428.889
$ws.Range("K89").Value = 7562.5
$ws.Range("L89").Value = 17144.445
$ws.Range("M89").Value = -1946.5
$ws.Range("N89").Value = -28376.445

$ws.Range("H107").Value = 1463.0256
$ws.Range("I107").Value = 823.8333
$ws.Range("J107").Value = 2485.7334
$ws.Range("K107").Value = 823.8333
$ws.Range("L107").Value = 2485.7334
$ws.Range("M107").Value = 1096.1667
$ws.Range("N107").Value = -6325.7334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 304
$ws.Range("I2").Value = 304
$ws.Range("K2").Value = 304
$ws.Range("M2").Value = -191

$ws.Range("H31").Value = 18382.06
$ws.Range("I31").Value = 1052.3077
$ws.Range("J31").Value = 43413.926
$ws.Range("K31").Value = 1052.3077
$ws.Range("L31").Value = 43413.926
$ws.Range("M31").Value = -757.3077000000001
$ws.Range("N31").Value = -44003.926

$ws.Range("H34").Value = 18382.06
$ws.Range("I34").Value = 1052.3077
$ws.Range("J34").Value = 43413.926
$ws.Range("K34").Value = 1052.3077
$ws.Range("L34").Value = 43413.926
$ws.Range("M34").Value = -850.3077000000001
$ws.Range("N34").Value = -43817.926

$ws.Range("H86").Value = 3968
$ws.Range("I86").Value = 3891.111
$ws.Range("K86").Value = 3891.111
$ws.Range("M86").Value = -2768.111

$ws.Range("H89").Value = 3968
$ws.Range("I89").Value = 3891.111
$ws.Range("K89").Value = 19455.555
$ws.Range("M89").Value = -13839.555

$ws.Range("H132").Value = 21280378
$ws.Range("I132").Value = 25003938
$ws.Range("J132").Value = 2889.4285
$ws.Range("K132").Value = 75011814
$ws.Range("L132").Value = 8668.2855
$ws.Range("M132").Value = -75009284
$ws.Range("N132").Value = -13728.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 700.3103599999999
$ws.Range("I131").Value = 452.56818
$ws.Range("J131").Value = 953.81396
$ws.Range("K131").Value = 1357.70454
$ws.Range("L131").Value = 2861.44188
$ws.Range("M131").Value = 3682.29546
$ws.Range("N131").Value = -12941.44188

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()

$ws.Range("H132").Value = 4963.8857
$ws.Range("I132").Value = 6194
$ws.Range("J132").Value = 2280
$ws.Range("K132").Value = 18582
$ws.Range("L132").Value = 6840
$ws.Range("M132").Value = -16052
$ws.Range("N132").Value = -11900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 615
$ws.Range("I16").Value = 538
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 538
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -368
$ws.Range("N16").Value = -1340

$ws.Range("H22").Value = 759
$ws.Range("I22").Value = 570
$ws.Range("J22").Value = 910.2
$ws.Range("K22").Value = 570
$ws.Range("L22").Value = 910.2
$ws.Range("M22").Value = -275
$ws.Range("N22").Value = -1500.2

$ws.Range("H27").Value = 759
$ws.Range("I27").Value = 570
$ws.Range("J27").Value = 910.2
$ws.Range("K27").Value = 570
$ws.Range("L27").Value = 910.2
$ws.Range("M27").Value = -463
$ws.Range("N27").Value = -1124.2

$ws.Range("H46").Value = 1825.125
$ws.Range("I46").Value = 2070.1667
$ws.Range("J46").Value = 1090
$ws.Range("K46").Value = 2070.1667
$ws.Range("L46").Value = 1090
$ws.Range("M46").Value = -1882.1667
$ws.Range("N46").Value = -1466

$ws.Range("H55").Value = 139.09091
$ws.Range("I55").Value = 66.666664
$ws.Range("J55").Value = 226
$ws.Range("K55").Value = 66.666664
$ws.Range("L55").Value = 226
$ws.Range("M55").Value = 106.333336
$ws.Range("N55").Value = -572

$ws.Range("H68").Value = 1950
$ws.Range("I68").Value = 1950
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1201
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1950
$ws.Range("I71").Value = 1950
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9750
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6006
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1135.3077
$ws.Range("I126").Value = 921
$ws.Range("J126").Value = 1230.5555
$ws.Range("K126").Value = 2763
$ws.Range("L126").Value = 3691.6665
$ws.Range("M126").Value = -293
$ws.Range("N126").Value = -8631.666499999999

$ws.Range("H131").Value = 37715
$ws.Range("J131").Value = 37715
$ws.Range("L131").Value = 37715
$ws.Range("N131").Value = -47795

$ws.Range("H132").Value = 2716.3125
$ws.Range("I132").Value = 3154.0217
$ws.Range("J132").Value = 1597.7222
$ws.Range("K132").Value = 9462.0651
$ws.Range("L132").Value = 4793.1666
$ws.Range("M132").Value = -6932.0651
$ws.Range("N132").Value = -9853.1666
